# Week 16 update: add rookie TE M.Hall as a newly-logged player.
# He's inserted as a new column between F.Darby and K.Pitts on both the
# "Rushing" and "Receiving" trackers, shifting K.Pitts/H.Hurst/L.Smith/
# P.Hesse/K.Smith one column to the right, with the row-2 "Yards list"
# marker left at "n" (not yet logged) like every other player.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Rushing", "Receiving")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # F.Darby is column O, K.Pitts is column P -> insert the new player's
    # column at P so it lands right after F.Darby.
    $ws.Columns("P").Insert()

    $ws.Range("P1").Value = "M.Hall"
    $ws.Range("P2").Value = "n"
}
